$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price / volume data (GitHub Actions refresh).
# D-column price values that look numeric must be written as literal text
# (matching the source data, e.g. "8.50" not 8.5) - force text format, set
# the value, then clear the temporary formatting so no style residue is left
# behind (cells keep their original default/unstyled appearance).

$ws.Range("D2").Value = '63.415.45'
$ws.Range("E2").Value = '  +0.74%  '

$ws.Range("D3").Value = '2.638.39'
$ws.Range("E3").Value = '  +2.10%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.59'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.28%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.84'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.66%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.589'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.99%  '

$ws.Range("D9").Value = '2.638.72'
$ws.Range("E9").Value = '  +2.14%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.107'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.78%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.65'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.09%  '

$ws.Range("E12").Value = '  +0.05%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.355'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.54%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.62'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.10%  '

$ws.Range("D15").Value = '3.118.78'
$ws.Range("E15").Value = '  +2.40%  '

$ws.Range("D16").Value = '63.381.74'
$ws.Range("E16").Value = '  +0.95%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000146'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.72%  '

$ws.Range("D18").Value = '2.662.51'
$ws.Range("E18").Value = '  +3.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.34'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '341.84'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.37'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.08%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.71'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.32%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.92'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.21%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.68'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +5.60%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.57'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +8.52%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.166'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '554.09'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +18.43%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.50'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.91%  '

$ws.Range("E30").Value = '  +0.21%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.82'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.39%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.98'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.92%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.79'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +11.67%  '

$ws.Range("D34").Value = '0.0₃0810'
$ws.Range("E34").Value = '  -1.48%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '175.17'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.86'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +7.75%  '

$ws.Range("E37").Value = '  +0.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.402'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.26%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.14'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.70%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.77'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.88%  '

$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '170.64'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +8.20%  '

$ws.Range("B42").Value = 'USDe'
$ws.Range("C42").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.01%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.42'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.58%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.75'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.06%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.95'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +3.86%  '

$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.629'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.78%  '

$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0553'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.20%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0240'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.75%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0957'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.78%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.85'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.31%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.71'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.90%  '
